# Update "想去人数" (F column) values for three worksheets to reflect
# newly generated output (gh-pages update at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 97
$ws1.Range("F4").Value = 612
$ws1.Range("F6").Value = 9222
$ws1.Range("F7").Value = 832
$ws1.Range("F8").Value = 333
$ws1.Range("F9").Value = 1181
$ws1.Range("F10").Value = 1072
$ws1.Range("F11").Value = 138
$ws1.Range("F12").Value = 59
$ws1.Range("F13").Value = 14
$ws1.Range("F14").Value = 253
$ws1.Range("F15").Value = 369
$ws1.Range("F16").Value = 79
$ws1.Range("F17").Value = 246
$ws1.Range("F18").Value = 1206

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 7

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 97
$ws4.Range("F5").Value = 7
$ws4.Range("F6").Value = 612
$ws4.Range("F8").Value = 9222
$ws4.Range("F9").Value = 832
$ws4.Range("F10").Value = 333
$ws4.Range("F11").Value = 1181
$ws4.Range("F12").Value = 1072
$ws4.Range("F13").Value = 138
$ws4.Range("F14").Value = 59
$ws4.Range("F15").Value = 14
$ws4.Range("F16").Value = 253
$ws4.Range("F17").Value = 369
$ws4.Range("F18").Value = 79
$ws4.Range("F19").Value = 246
$ws4.Range("F20").Value = 1206
